$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 38-50 (ticker list trimmed from 48 to 36 entries)
$ws.Rows("38:50").Delete()

# Update ticker symbols for rows 2-37
$ws.Range("B2").Value = "NSE:ALPHAETF"
$ws.Range("C2").Value = "NSE:AKSHAR"
$ws.Range("D2").Value = ""
$ws.Range("F2").Value = "NSE:AMBUJACEM"
$ws.Range("B3").Value = "NSE:AMBUJACEM"
$ws.Range("C3").Value = "NSE:AVTNPL"
$ws.Range("E3").Value = "NSE:HINDALCO"
$ws.Range("B4").Value = "NSE:BAJAJ-AUTO"
$ws.Range("C4").Value = "NSE:ICICIPRULI"
$ws.Range("E4").Value = "NSE:HINDPETRO"
$ws.Range("F4").Value = "NSE:HINDUNILVR"
$ws.Range("B5").Value = "NSE:BANSWRAS"
$ws.Range("C5").Value = "NSE:JMFINANCIL"
$ws.Range("E5").Value = "NSE:IOC"
$ws.Range("F5").Value = ""
$ws.Range("B6").Value = "NSE:CARYSIL"
$ws.Range("C6").Value = "NSE:KECL"
$ws.Range("E6").Value = "NSE:LTIM"
$ws.Range("F6").Value = ""
$ws.Range("B7").Value = "NSE:CHEMBOND"
$ws.Range("C7").Value = "NSE:NAVINIFTY"
$ws.Range("E7").Value = "NSE:NESTLEIND"
$ws.Range("F7").Value = ""
$ws.Range("B8").Value = "NSE:CHENNPETRO"
$ws.Range("E8").Value = "NSE:OFSS"
$ws.Range("B9").Value = "NSE:DEVYANI"
$ws.Range("B10").Value = "NSE:DOLATALGO"
$ws.Range("B11").Value = "NSE:FILATEX"
$ws.Range("B12").Value = "NSE:FINOPB"
$ws.Range("B13").Value = "NSE:FIVESTAR"
$ws.Range("B14").Value = "NSE:GNA"
$ws.Range("B15").Value = "NSE:HARSHA"
$ws.Range("B16").Value = "NSE:HDFCNIFTY"
$ws.Range("B17").Value = "NSE:HINDUNILVR"
$ws.Range("B18").Value = "NSE:INDORAMA"
$ws.Range("B19").Value = "NSE:IONEXCHANG"
$ws.Range("B20").Value = "NSE:JHS"
$ws.Range("B21").Value = "NSE:JISLDVREQS"
$ws.Range("B22").Value = "NSE:KOKUYOCMLN"
$ws.Range("B23").Value = "NSE:KPIL"
$ws.Range("B24").Value = "NSE:KRBL"
$ws.Range("B25").Value = "NSE:MINDACORP"
$ws.Range("B26").Value = "NSE:MSPL"
$ws.Range("B27").Value = "NSE:NAVKARCORP"
$ws.Range("B28").Value = "NSE:NILAINFRA"
$ws.Range("B29").Value = "NSE:NV20BEES"
$ws.Range("B30").Value = "NSE:ORIENTALTL"
$ws.Range("B31").Value = "NSE:OSWALSEEDS"
$ws.Range("B32").Value = "NSE:RADICO"
$ws.Range("B33").Value = "NSE:RAILTEL"
$ws.Range("B34").Value = "NSE:RAYMOND"
$ws.Range("B35").Value = "NSE:RICOAUTO"
$ws.Range("B36").Value = "NSE:RUCHIRA"
$ws.Range("B37").Value = "NSE:SAH"
